# This script updates the "想去人数" (want-to-go count, column F) figures
# across the four worksheets of the 上海-漫展信息 workbook to reflect a
# newly generated data snapshot (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

function Get-Sheet($name) {
    return $wb.Worksheets.Item($name)
}

# 展览 (Exhibitions) sheet
$ws = Get-Sheet "展览"
$ws.Range("F13").Value = 110
$ws.Range("F14").Value = 7481
$ws.Range("F16").Value = 7717
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 57844
$ws.Range("F20").Value = 4812
$ws.Range("F23").Value = 505

# 演出 (Performances) sheet
$ws = Get-Sheet "演出"
$ws.Range("F9").Value = 7649
$ws.Range("F32").Value = 5
$ws.Range("F45").Value = 32

# 本地生活 (Local life) sheet
$ws = Get-Sheet "本地生活"
$ws.Range("F9").Value = 9445
$ws.Range("F15").Value = 275
$ws.Range("F16").Value = 2406
$ws.Range("F17").Value = 100
$ws.Range("F19").Value = 528

# 全部类型 (All types) sheet
$ws = Get-Sheet "全部类型"
$ws.Range("F8").Value = 275
$ws.Range("F11").Value = 110
$ws.Range("F12").Value = 7481
$ws.Range("F13").Value = 7717
$ws.Range("F14").Value = 57844
$ws.Range("F17").Value = 4812
$ws.Range("F20").Value = 505
$ws.Range("F29").Value = 528
$ws.Range("F50").Value = 32

$wb.Save()
